$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value2 = 74375.664
$ws.Range("I103").Value2 = 400000
$ws.Range("J103").Value2 = 9250.799999999999
$ws.Range("K103").Value2 = 1200000
$ws.Range("L103").Value2 = 27752.4
$ws.Range("M103").Value2 = -1199414
$ws.Range("N103").Value2 = -28924.4
$ws.Range("H121").Value2 = 1544.6666
$ws.Range("J121").Value2 = 1544.6666
$ws.Range("L121").Value2 = 4633.9998
$ws.Range("N121").Value2 = -8127.9998
$ws.Range("H137").Value2 = 243258.56
$ws.Range("J137").Value2 = 8871.6
$ws.Range("L137").Value2 = 26614.8
$ws.Range("N137").Value2 = -31714.8
$ws.Range("H141").Value2 = 6359.2173
$ws.Range("I141").Value2 = 4418.7334
$ws.Range("K141").Value2 = 13256.2002
$ws.Range("M141").Value2 = -8076.200199999999

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value2 = 35841.594
$ws.Range("I2").Value2 = 6162.8
$ws.Range("J2").Value2 = 85306.25
$ws.Range("K2").Value2 = 6162.8
$ws.Range("L2").Value2 = 85306.25
$ws.Range("M2").Value2 = -6049.8
$ws.Range("N2").Value2 = -85532.25
$ws.Range("H32").Value2 = 2573.6956
$ws.Range("I32").Value2 = 2499.9302
$ws.Range("K32").Value2 = 2499.9302
$ws.Range("M32").Value2 = -2212.9302
$ws.Range("H45").Value2 = 9779.762000000001
$ws.Range("I45").Value2 = 17199.6
$ws.Range("J45").Value2 = 3034.4546
$ws.Range("K45").Value2 = 17199.6
$ws.Range("L45").Value2 = 3034.4546
$ws.Range("M45").Value2 = -16822.6
$ws.Range("N45").Value2 = -3788.4546
$ws.Range("H61").Value2 = 8247.862999999999
$ws.Range("I61").Value2 = 10148.143
$ws.Range("K61").Value2 = 10148.143
$ws.Range("M61").Value2 = -9936.143
$ws.Range("H74").Value2 = 112289.92
$ws.Range("I74").Value2 = 121230.75
$ws.Range("K74").Value2 = 121230.75
$ws.Range("M74").Value2 = -120356.75
$ws.Range("H77").Value2 = 112289.92
$ws.Range("I77").Value2 = 121230.75
$ws.Range("K77").Value2 = 606153.75
$ws.Range("M77").Value2 = -601785.75
$ws.Range("H97").Value2 = 5151.4346
$ws.Range("I97").Value2 = 5835.9473
$ws.Range("K97").Value2 = 5835.9473
$ws.Range("M97").Value2 = -5339.9473
$ws.Range("H116").Value2 = 35841.594
$ws.Range("I116").Value2 = 6162.8
$ws.Range("J116").Value2 = 85306.25
$ws.Range("K116").Value2 = 6162.8
$ws.Range("L116").Value2 = 85306.25
$ws.Range("M116").Value2 = -3868.8
$ws.Range("N116").Value2 = -89894.25
$ws.Range("H122").Value2 = 741717.9
$ws.Range("I122").Value2 = 3932.7778
$ws.Range("J122").Value2 = 1405724.5
$ws.Range("K122").Value2 = 11798.3334
$ws.Range("L122").Value2 = 4217173.5
$ws.Range("M122").Value2 = -9348.3334
$ws.Range("N122").Value2 = -4222073.5
$ws.Range("H136").Value2 = 8247.862999999999
$ws.Range("I136").Value2 = 10148.143
$ws.Range("K136").Value2 = 30444.429
$ws.Range("M136").Value2 = -27894.429

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value2 = 35841.594
$ws.Range("I3").Value2 = 6162.8
$ws.Range("J3").Value2 = 85306.25
$ws.Range("K3").Value2 = 6162.8
$ws.Range("L3").Value2 = 85306.25
$ws.Range("M3").Value2 = -6048.8
$ws.Range("N3").Value2 = -85534.25
$ws.Range("H54").Value2 = 18193.666
$ws.Range("I54").Value2 = 15832.6
$ws.Range("K54").Value2 = 15832.6
$ws.Range("M54").Value2 = -15348.6
$ws.Range("H134").Value2 = 11800.3125
$ws.Range("I134").Value2 = 13195.462
$ws.Range("J134").Value2 = 5754.6665
$ws.Range("K134").Value2 = 39586.386
$ws.Range("L134").Value2 = 17263.9995
$ws.Range("M134").Value2 = -37051.386
$ws.Range("N134").Value2 = -22333.9995
$ws.Range("H140").Value2 = 66000
$ws.Range("I140").Value2 = 45000
$ws.Range("J140").Value2 = 87000
$ws.Range("K140").Value2 = 45000
$ws.Range("L140").Value2 = 87000
$ws.Range("M140").Value2 = -39820
$ws.Range("N140").Value2 = -97360
$ws.Range("H141").Value2 = 100375
$ws.Range("J141").Value2 = 100375
$ws.Range("L141").Value2 = 100375
$ws.Range("N141").Value2 = -110735

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value2 = 4501
$ws.Range("I3").Value2 = 4501
$ws.Range("K3").Value2 = 4501
$ws.Range("M3").Value2 = -4388
$ws.Range("H31").Value2 = 2102.7656
$ws.Range("I31").Value2 = 1383.0588
$ws.Range("J31").Value2 = 2918.4333
$ws.Range("K31").Value2 = 1383.0588
$ws.Range("L31").Value2 = 2918.4333
$ws.Range("M31").Value2 = -1088.0588
$ws.Range("N31").Value2 = -3508.4333
$ws.Range("H34").Value2 = 2102.7656
$ws.Range("I34").Value2 = 1383.0588
$ws.Range("J34").Value2 = 2918.4333
$ws.Range("K34").Value2 = 1383.0588
$ws.Range("L34").Value2 = 2918.4333
$ws.Range("M34").Value2 = -1181.0588
$ws.Range("N34").Value2 = -3322.4333
$ws.Range("H107").Value2 = 55568348
$ws.Range("J107").Value2 = 1649.75
$ws.Range("L107").Value2 = 1649.75
$ws.Range("N107").Value2 = -5489.75
$ws.Range("H132").Value2 = 27907.65
$ws.Range("I132").Value2 = 9877.066000000001
$ws.Range("K132").Value2 = 29631.198
$ws.Range("M132").Value2 = -27101.198

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value2 = 313551.03
$ws.Range("J5").Value2 = 835174.5
$ws.Range("L5").Value2 = 2505523.5
$ws.Range("N5").Value2 = -2505747.5
$ws.Range("H68").Value2 = 17247594
$ws.Range("J68").Value2 = 20840454
$ws.Range("L68").Value2 = 62521362
$ws.Range("N68").Value2 = -62522984
$ws.Range("H71").Value2 = 17247594
$ws.Range("J71").Value2 = 20840454
$ws.Range("L71").Value2 = 187564086
$ws.Range("N71").Value2 = -187572198
$ws.Range("H98").Value2 = 992.6
$ws.Range("J98").Value2 = 864.8
$ws.Range("L98").Value2 = 2594.4
$ws.Range("N98").Value2 = -5590.4
$ws.Range("H107").Value2 = 1479.4736
$ws.Range("I107").Value2 = 500
$ws.Range("J107").Value2 = 1740.6666
$ws.Range("K107").Value2 = 1500
$ws.Range("L107").Value2 = 5221.9998
$ws.Range("M107").Value2 = 420
$ws.Range("N107").Value2 = -9061.9998
$ws.Range("H113").Value2 = 829.45715
$ws.Range("I113").Value2 = 582.1667
$ws.Range("J113").Value2 = 958.4783
$ws.Range("K113").Value2 = 1746.5001
$ws.Range("L113").Value2 = 2875.4349
$ws.Range("M113").Value2 = 423.4999
$ws.Range("N113").Value2 = -7215.4349
$ws.Range("H133").Value2 = 3909.3635
$ws.Range("I133").Value2 = 3909.3635
$ws.Range("J133").Value2 = 0
$ws.Range("K133").Value2 = 11728.0905
$ws.Range("L133").Value2 = 0
$ws.Range("M133").Value2 = -6668.0905
$ws.Range("N133").ClearContents()
$ws.Range("H135").Value2 = 313551.03
$ws.Range("J135").Value2 = 835174.5
$ws.Range("L135").Value2 = 7516570.5
$ws.Range("N135").Value2 = -7521640.5
$ws.Range("H137").Value2 = 6650.4136
$ws.Range("J137").Value2 = 10011.125
$ws.Range("L137").Value2 = 30033.375
$ws.Range("N137").Value2 = -40233.375

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value2 = 906.93335
$ws.Range("I5").Value2 = 302
$ws.Range("K5").Value2 = 302
$ws.Range("M5").Value2 = -190
$ws.Range("H43").Value2 = 12067.583
$ws.Range("I43").Value2 = 10545.286
$ws.Range("J43").Value2 = 14198.8
$ws.Range("K43").Value2 = 10545.286
$ws.Range("L43").Value2 = 14198.8
$ws.Range("M43").Value2 = -10394.286
$ws.Range("N43").Value2 = -14500.8
$ws.Range("H57").Value2 = 49999
$ws.Range("J57").Value2 = 49999
$ws.Range("L57").Value2 = 49999
$ws.Range("N57").Value2 = -51639
$ws.Range("H122").Value2 = 12177.182
$ws.Range("I122").Value2 = 12376.143
$ws.Range("K122").Value2 = 37128.429
$ws.Range("M122").Value2 = -34678.429

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value2 = 1348.8077
$ws.Range("I61").Value2 = 1074.3889
$ws.Range("K61").Value2 = 1074.3889
$ws.Range("M61").Value2 = -872.3888999999999
$ws.Range("H113").Value2 = 1348.8077
$ws.Range("I113").Value2 = 1074.3889
$ws.Range("K113").Value2 = 1074.3889
$ws.Range("M113").Value2 = 1095.6111
$ws.Range("H122").Value2 = 5012
$ws.Range("I122").Value2 = 3159.6
$ws.Range("J122").Value2 = 5854
$ws.Range("K122").Value2 = 9478.799999999999
$ws.Range("L122").Value2 = 17562
$ws.Range("M122").Value2 = -7028.799999999999
$ws.Range("N122").Value2 = -22462
$ws.Range("H141").Value2 = 48998
$ws.Range("J141").Value2 = 48998
$ws.Range("L141").Value2 = 48998
$ws.Range("N141").Value2 = -59358

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value2 = 71436250
$ws.Range("I2").Value2 = 83342040
$ws.Range("K2").Value2 = 83342040
$ws.Range("M2").Value2 = -83341928
$ws.Range("H4").Value2 = 4686.6313
$ws.Range("I4").Value2 = 5390
$ws.Range("J4").Value2 = 2049
$ws.Range("K4").Value2 = 5390
$ws.Range("L4").Value2 = 2049
$ws.Range("M4").Value2 = -5277
$ws.Range("N4").Value2 = -2275
$ws.Range("H29").Value2 = 16895.445
$ws.Range("I29").Value2 = 13437.143
$ws.Range("J29").Value2 = 28999.5
$ws.Range("K29").Value2 = 13437.143
$ws.Range("L29").Value2 = 28999.5
$ws.Range("M29").Value2 = -13147.143
$ws.Range("N29").Value2 = -29579.5
$ws.Range("H93").Value2 = 72500
$ws.Range("J93").Value2 = 72500
$ws.Range("L93").Value2 = 72500
$ws.Range("N93").Value2 = -77492
$ws.Range("H131").Value2 = 25000
$ws.Range("J131").Value2 = 25000
$ws.Range("L131").Value2 = 25000
$ws.Range("N131").Value2 = -35080
$ws.Range("H132").Value2 = 17624.346
$ws.Range("J132").Value2 = 6553.8887
$ws.Range("L132").Value2 = 19661.6661
$ws.Range("N132").Value2 = -24721.6661
$ws.Range("H136").Value2 = 3361.2888
$ws.Range("I136").Value2 = 2155.182
$ws.Range("J136").Value2 = 6678.0835
$ws.Range("K136").Value2 = 6465.545999999999
$ws.Range("L136").Value2 = 20034.2505
$ws.Range("M136").Value2 = -3915.545999999999
$ws.Range("N136").Value2 = -25134.2505
